$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.211.93'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '1.807.66'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '''224.34'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '''31.98'
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('D9').Value = '''0.291'
$ws.Range('E9').Value = '  +2.50%  '
$ws.Range('D10').Value = '''0.0718'
$ws.Range('E10').Value = '  +8.25%  '
$ws.Range('D11').Value = '''0.0928'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').Value = '2.070.27'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').Value = '1.808.19'
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('D14').Value = '''10.83'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '''0.634'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '34.214.26'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '''69.32'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').Value = '''248.39'
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('D20').Value = '0.0₃0796'
$ws.Range('E20').Value = '  +6.67%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '''0.998'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = '''10.92'
$ws.Range('E22').Value = '  +4.68%  '
$ws.Range('D23').Value = '''4.22'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('D24').Value = '''2.14'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = '''159.37'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('D26').Value = '''16.60'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').Value = '''7.16'
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').Value = '''0.0528'
$ws.Range('E30').Value = '  +2.68%  '
$ws.Range('D31').Value = '''3.76'
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('D33').Value = '''3.57'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').Value = '1.424.61'
$ws.Range('E35').Value = '  -2.07%  '
$ws.Range('E36').Value = '  +0.65%  '
$ws.Range('D37').Value = '''0.635'
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').Value = '''0.952'
$ws.Range('E39').Value = '  +6.40%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '''80.85'
$ws.Range('E40').Value = '  -2.69%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '''2.75'
$ws.Range('E41').Value = '  -3.82%  '
$ws.Range('D42').Value = '''2.34'
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('E43').Value = '  +3.75%  '
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('D45').Value = '''0.0498'
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('D47').Value = '1.966.25'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').Value = '''106.26'
$ws.Range('E48').Value = '  +6.75%  '
$ws.Range('D49').Value = '''0.996'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').Value = '''11.85'
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('E51').Value = '  +5.62%  '
